$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.117.69"
$ws.Range("E2").Value = "'  -1.93%  "
$ws.Range("D3").Value = "'1.835.00"
$ws.Range("E3").Value = "'  -3.35%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "'  +0.08%  "
$ws.Range("D5").Value = "'228.17"
$ws.Range("E5").Value = "'  -4.57%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "'  +0.17%  "
$ws.Range("D7").Value = "'0.4598"
$ws.Range("E7").Value = "'  -4.85%  "
$ws.Range("D8").Value = "'0.2676"
$ws.Range("E8").Value = "'  -6.32%  "
$ws.Range("D9").Value = "'0.06208"
$ws.Range("E9").Value = "'  -5.38%  "
$ws.Range("D10").Value = "'1.838.93"
$ws.Range("E10").Value = "'  -4.93%  "
$ws.Range("D11").Value = "'0.07333"
$ws.Range("E11").Value = "'  -1.95%  "
$ws.Range("D12").Value = "'15.90"
$ws.Range("E12").Value = "'  -5.14%  "
$ws.Range("D13").Value = "'4.873"
$ws.Range("E13").Value = "'  -4.80%  "
$ws.Range("D14").Value = "'82.81"
$ws.Range("E14").Value = "'  -6.13%  "
$ws.Range("D15").Value = "'0.6149"
$ws.Range("E15").Value = "'  -8.08%  "
$ws.Range("D16").Value = "'30.038.53"
$ws.Range("E16").Value = "'  -2.18%  "
$ws.Range("E17").Value = "'  +0.23%  "
$ws.Range("D18").Value = "'225.33"
$ws.Range("E18").Value = "'  -2.62%  "
$ws.Range("B19").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D19").Value = "'2.084.73"
$ws.Range("E19").Value = "'  -6.88%  "
$ws.Range("B20").Value = "'BinanceUSD"
$ws.Range("C20").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "'  +0.01%  "
$ws.Range("D21").Value = "'12.28"
$ws.Range("E21").Value = "'  -7.90%  "
$ws.Range("D22").Value = "'0.000007190"
$ws.Range("E22").Value = "'  -5.57%  "
$ws.Range("D23").Value = "'4.834"
$ws.Range("E23").Value = "'  -8.71%  "
$ws.Range("B24").Value = "'BitDAO"
$ws.Range("C24").Value = "'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D24").Value = "'0.3925"
$ws.Range("E24").Value = "'  +2.17%  "
$ws.Range("B25").Value = "'Chainlink"
$ws.Range("C25").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").Value = "'5.827"
$ws.Range("E25").Value = "'  -6.66%  "
$ws.Range("B26").Value = "'Monero"
$ws.Range("C26").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'165.12"
$ws.Range("E26").Value = "'  -2.86%  "
$ws.Range("B27").Value = "'Cosmos"
$ws.Range("C27").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'9.064"
$ws.Range("E27").Value = "'  -3.14%  "
$ws.Range("B28").Value = "'EthereumClassic"
$ws.Range("C28").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'17.54"
$ws.Range("E28").Value = "'  -6.72%  "
$ws.Range("B29").Value = "'LidoDAOToken"
$ws.Range("C29").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'1.842"
$ws.Range("E29").Value = "'  -6.41%  "
$ws.Range("B30").Value = "'Stellar"
$ws.Range("C30").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.1009"
$ws.Range("E30").Value = "'  -1.07%  "
$ws.Range("B31").Value = "'Toncoin"
$ws.Range("C31").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'1.375"
$ws.Range("E31").Value = "'  -2.02%  "
$ws.Range("B32").Value = "'InternetComputer(DFINITY)"
$ws.Range("C32").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'4.049"
$ws.Range("E32").Value = "'  -7.11%  "
$ws.Range("B33").Value = "'Filecoin"
$ws.Range("C33").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'3.739"
$ws.Range("E33").Value = "'  -7.35%  "
$ws.Range("B34").Value = "'Hedera"
$ws.Range("C34").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.04790"
$ws.Range("E34").Value = "'  -6.49%  "
$ws.Range("B35").Value = "'ARBITRUM"
$ws.Range("C35").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.131"
$ws.Range("E35").Value = "'  -7.46%  "
$ws.Range("D36").Value = "'0.6899"
$ws.Range("E36").Value = "'  -9.25%  "
$ws.Range("B37").Value = "'Frax"
$ws.Range("C37").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "'1.002"
$ws.Range("E37").Value = "'  -0.18%  "
$ws.Range("B38").Value = "'HuobiToken"
$ws.Range("C38").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").Value = "'2.703"
$ws.Range("E38").Value = "'  -0.29%  "
$ws.Range("B39").Value = "'VeChain"
$ws.Range("C39").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01807"
$ws.Range("E39").Value = "'  -4.18%  "
$ws.Range("B40").Value = "'MXToken"
$ws.Range("C40").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.597"
$ws.Range("E40").Value = "'  -2.36%  "
$ws.Range("B41").Value = "'TrustWalletToken"
$ws.Range("C41").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.8805"
$ws.Range("E41").Value = "'  -4.40%  "
$ws.Range("D42").Value = "'1.910"
$ws.Range("E42").Value = "'  -8.28%  "
$ws.Range("B43").Value = "'PaxDollar"
$ws.Range("C43").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "'0.9969"
$ws.Range("E43").Value = "'  -0.80%  "
$ws.Range("B44").Value = "'Quant"
$ws.Range("C44").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'102.85"
$ws.Range("E44").Value = "'  -4.25%  "
$ws.Range("B45").Value = "'FraxShare"
$ws.Range("C45").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'5.432"
$ws.Range("E45").Value = "'  -5.04%  "
$ws.Range("B46").Value = "'TheSandbox"
$ws.Range("C46").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").Value = "'0.3972"
$ws.Range("E46").Value = "'  -7.88%  "
$ws.Range("B47").Value = "'Aptos"
$ws.Range("C47").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'6.850"
$ws.Range("E47").Value = "'  -7.94%  "
$ws.Range("B48").Value = "'Algorand"
$ws.Range("C48").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.1176"
$ws.Range("E48").Value = "'  -7.86%  "
$ws.Range("B49").Value = "'Aave"
$ws.Range("C49").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'58.88"
$ws.Range("E49").Value = "'  -9.21%  "
$ws.Range("B50").Value = "'EnergySwap"
$ws.Range("C50").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.403"
$ws.Range("E50").Value = "'  -6.46%  "
$ws.Range("B51").Value = "'Cronos"
$ws.Range("C51").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05532"
$ws.Range("E51").Value = "'  -2.58%  "
